$wb = $excel.ActiveWorkbook

# ---- Sheet: Summary ----
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.7419928825622776
$ws1.Range("C2").Value = 0.1320754716981132
$ws1.Range("D2").Value = 0.75
$ws1.Range("E2").Value = 0.2245989304812834
$ws1.Range("F2").Value = 0.3874538745387454
$ws1.Range("G2").Value = 0.6356228172293364
$ws1.Range("H2").Value = 0.7978196896736223
$ws1.Range("I2").Value = 21
$ws1.Range("J2").Value = 138
$ws1.Range("K2").Value = 396
$ws1.Range("L2").Value = 7

# ---- Sheet: Classification Report ----
$ws2 = $wb.Worksheets.Item("Classification Report")

# row 2 - label "0"
$ws2.Range("B2").Value = 0.9826302729528535
$ws2.Range("C2").Value = 0.7415730337078652
$ws2.Range("D2").Value = 0.8452508004268944

# row 3 - label "1"
$ws2.Range("B3").Value = 0.1320754716981132
$ws2.Range("C3").Value = 0.75
$ws2.Range("D3").Value = 0.2245989304812834

# row 4 - accuracy
$ws2.Range("B4").Value = 0.7419928825622776
$ws2.Range("C4").Value = 0.7419928825622776
$ws2.Range("D4").Value = 0.7419928825622776
$ws2.Range("E4").Value = 0.7419928825622776

# row 5 - macro avg
$ws2.Range("B5").Value = 0.5573528723254834
$ws2.Range("C5").Value = 0.7457865168539326
$ws2.Range("D5").Value = 0.5349248654540889

# row 6 - weighted avg
$ws2.Range("B6").Value = 0.9402538771608026
$ws2.Range("C6").Value = 0.7419928825622776
$ws2.Range("D6").Value = 0.8143286432054049

# ---- Sheet: Confusion Matrix ----
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 396
$ws3.Range("C2").Value = 138
$ws3.Range("B3").Value = 7
$ws3.Range("C3").Value = 21
